# Adding the changes we made on may 9th
#
# Insert 17 fresh accelerometer samples at the top of the data block
# (pushing the previously-first 13 rows down), and drop the final 7 rows
# of the old trailing (near-zero) tail so the sheet ends up with 30 data
# rows (A1:C31) instead of the original 20 (A1:C21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to place at the top of the data (becomes rows 2..18)
$newData = @(
    @(5.212075734138489, -2.939898788928986, 1.854160755872726),
    @(3.245875406265258, -4.094250345230103, 2.58136396408081),
    @(3.096780717372894, -3.371150910854339, 2.613386332988739),
    @(2.924881196022033, -2.50173692703247, 2.383840799331665),
    @(3.458236134052277, -2.494953083992004, 2.681476718187333),
    @(2.869051647186279, -3.347217082977296, 2.930787801742554),
    @(2.807376968860626, -3.875039219856262, 3.420289939641953),
    @(2.295876741409302, -4.034408569335938, 3.253981232643127),
    @(2.279258108139039, -3.775099605321884, 3.11082683801651),
    @(2.566667938232422, -3.378203916549682, 3.007539582252503),
    @(3.106618106365205, -3.249815458059311, 3.031012719869614),
    @(2.987140679359436, -3.142817544937134, 3.183629143238068),
    @(2.434188187122345, -3.181812554597855, 3.162444919347763),
    @(2.282221984863281, -3.265003252029419, 3.094355344772339),
    @(2.110153055191039, -3.195758980512619, 3.138975620269776),
    @(1.555334329605102, -2.938729083538055, 3.47747951745987),
    @(0.6493126988410929, -2.875420850515366, 3.464587104320525)
)

$insertCount = $newData.Length

# Insert blank rows right after the header (before the current row 2),
# shifting all existing data rows down by $insertCount.
$insertEndRow = 1 + $insertCount
$insertRange = $ws.Range("A2:C" + $insertEndRow)
$insertRange.EntireRow.Insert()
# The insert operation copies the header row's formatting (bold/border)
# onto the newly-created rows; clear that so the new data rows look like
# the plain, unstyled data rows elsewhere in the sheet.
$insertRange.ClearFormats()

# Fill the newly inserted rows with the fresh data.
for ($i = 0; $i -lt $insertCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $newData[$i][0]
    $ws.Cells.Item($r, 2).Value = $newData[$i][1]
    $ws.Cells.Item($r, 3).Value = $newData[$i][2]
}

# The old tail (previously rows 15-21, now shifted down by $insertCount
# to rows 32-38) is no longer part of the dataset - remove it.
$oldTailStart = 15 + $insertCount
$oldTailEnd = 21 + $insertCount
$ws.Range("A" + $oldTailStart + ":C" + $oldTailEnd).EntireRow.Delete()
